$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price-column (D) writes stay text, matching the source data
# (avoids Excel auto-coercing strings like "1.0000" or "234.77" into numbers)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.897.01"
$ws.Range("E2").Value = "  -4.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.631.77"
$ws.Range("E3").Value = "  -6.47%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.77"
$ws.Range("E5").Value = "  -5.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4729"
$ws.Range("E7").Value = "  -6.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2542"
$ws.Range("E8").Value = "  -7.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06117"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06953"
$ws.Range("E10").Value = "  -4.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.634.01"
$ws.Range("E11").Value = "  -6.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.65"
$ws.Range("E12").Value = "  -3.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6110"
$ws.Range("E13").Value = "  -6.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.334"
$ws.Range("E14").Value = "  -6.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "72.59"
$ws.Range("E15").Value = "  -6.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9980"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.893.69"
$ws.Range("E18").Value = "  -4.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006560"
$ws.Range("E19").Value = "  -4.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.06"
$ws.Range("E20").Value = "  -6.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.840.92"
$ws.Range("E21").Value = "  -6.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.333"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.539"
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.225"
$ws.Range("E24").Value = "  -3.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "134.19"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.73"
$ws.Range("E26").Value = "  -3.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.369"
$ws.Range("E27").Value = "  -8.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "102.60"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.632"
$ws.Range("E29").Value = "  -8.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.743"
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07712"
$ws.Range("E31").Value = "  -6.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.530"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9990"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04271"
$ws.Range("E34").Value = "  -8.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.599"
$ws.Range("E35").Value = "  -2.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9193"
$ws.Range("E36").Value = "  -7.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5765"
$ws.Range("E37").Value = "  -6.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.537"
$ws.Range("E38").Value = "  -7.87%  "
$ws.Range("E39").Value = "  -4.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9977"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8190"
$ws.Range("E41").Value = "  +8.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.30"
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.774"
$ws.Range("E43").Value = "  -8.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3679"
$ws.Range("E44").Value = "  -6.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.693"
$ws.Range("E45").Value = "  -6.30%  "
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05193"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.012"
$ws.Range("E48").Value = "  -4.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.39"
$ws.Range("E49").Value = "  -4.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("E51").Value = "  -0.13%  "
